# Monopoly Datasheet edit: insert a new "Location" column (D) on the
# "Card" sheet, shifting the existing Cost-per-House / Base / House1-4 /
# Hotel / Colour columns (D:K) one column to the right (E:L), then fill
# in the new Location column with each property's board position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card")
$ws.Activate()

# --- Shift columns D:K one column right, into E:L, for every data row.
# Work right-to-left (column K -> L first) so we never clobber a value
# before it has been read, and work cell-by-cell (Value2) because a
# block Range.Value round-trip isn't reliable in this host. Also carry
# the bold flag along with the value (row 9's Hotel price is bold) and
# clear it from the vacated source cell.
for ($r = 1; $r -le 29; $r++) {
    for ($c = 11; $c -ge 4; $c--) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($r, $c + 1)
        $wasBold = $srcCell.Font.Bold
        $dstCell.Value2 = $srcCell.Value2
        $dstCell.Font.Bold = $wasBold
        $srcCell.Font.Bold = $false
    }
}

# --- New column header
$ws.Range("D1").Value = "Location"

# --- New column values: each property's position on the Monopoly board
$locations = @{
    2 = 1;  3 = 3;  4 = 6;  5 = 8;  6 = 9;  7 = 11; 8 = 13; 9 = 14;
    10 = 16; 11 = 18; 12 = 19; 13 = 21; 14 = 23; 15 = 24; 16 = 26;
    17 = 27; 18 = 29; 19 = 31; 20 = 32; 21 = 34; 22 = 37; 23 = 39;
    24 = 5; 25 = 15; 26 = 25; 27 = 35; 28 = 12; 29 = 28
}
foreach ($r in $locations.Keys) {
    $ws.Cells.Item($r, 4).Value2 = $locations[$r]
}

# --- View: normal zoom 100%, selection on D4
$excel.ActiveWindow.Zoom = 100
$ws.Range("D4").Select()
